$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (new row), entered first (E16 then D16) so shared-string
# indices line up the same way they did in the authored workbook:
# E16 -> "master_all_responses_SB_resub_Oct-01-2023.csv"
# D16 -> "re_submitted_tracker_SB_Oct-01-2023.csv"
$ws.Range("E16").Value = "master_all_responses_SB_resub_Oct-01-2023.csv"
$ws.Range("D16").Value = "re_submitted_tracker_SB_Oct-01-2023.csv"

# --- Row 15 (existing row), add new response_collected value
# E15 -> "master_worker_response_tracke_SB_Oct-01-2023.csv"
$ws.Range("E15").Value = "master_worker_response_tracke_SB_Oct-01-2023.csv"

# --- Row 17 (new row)
# D17 -> "all_submitted_tracker_SB_video_no_play_Oct-01-2023.csv"
$ws.Range("D17").Value = "all_submitted_tracker_SB_video_no_play_Oct-01-2023.csv"

# Remaining row 16 cells (reuse existing shared strings)
$ws.Range("A16").Value = "paiewise_resub"
$ws.Range("C16").Value = "SB"

# Remaining row 17 cells (reuse existing shared strings)
$ws.Range("A17").Value = "pairwise_resub_video_not_play"
$ws.Range("C17").Value = "SB"

# Date-like text values ("Oct-1-2023") must stay text, not get parsed
# into a serial date. A leading apostrophe forces text entry (same as
# typing it straight into Excel); reset the style back to Normal
# afterwards so no stray style index is left on the cell (matches the
# plain <c r="B16" t="s"> cell in the target).
$ws.Range("B16").Value = "'Oct-1-2023"
$ws.Range("B16").Style = "Normal"

$ws.Range("B17").Value = "'Oct-1-2023"
$ws.Range("B17").Style = "Normal"

# Page setup + selection metadata
$ws.PageSetup.Orientation = 1

$ws.Range("E18").Select()
